$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at G (pushes the existing "Sensor Link" column to H) ---
$ws.Columns.Item(7).Insert()

# --- New column G header + content: "Which connectors/signals" used in the IMU section ---
$ws.Range("G1").Value = "Which connectors/signals"

$ws.Range("G12").Value = ":)"
$ws.Range("G13").Value = ":)"
$ws.Range("G14").Value = ":) connect to 11"
$ws.Range("G15").Value = ":) connect to 13"
$ws.Range("G21").Value = ":) connect to 9"
$ws.Range("G22").Value = ":) connect to 10"
$ws.Range("G23").Value = ":) connect to 12"
$ws.Range("G24").Value = "bridge to SDOM"

# Match the wrap-text / bordered style used by the neighboring F23/F24 cells
$ws.Range("G23").Style = $ws.Range("F23").Style
$ws.Range("G24").Style = $ws.Range("F24").Style

# --- Rename the IMU "VDD" signal to "VDD/VCC" ---
$ws.Range("E13").Value = "VDD/VCC"

# --- Turn the two part-link URLs into real hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("D6"), "http://www.digikey.com/scripts/DkSearch/dksus.dll?Detail&itemSeq=252199103&uq=636543209955835002")
$ws.Hyperlinks.Add($ws.Range("H11"), "https://www.sparkfun.com/products/13284")
